$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.855.59'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '1.869.03'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.66'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5079'
$ws.Range("E7").Value = '  -1.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3664'
$ws.Range("E8").Value = '  -2.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07175'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8906'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.66'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07505'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.878.01'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.93'
$ws.Range("E14").Value = '  +5.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.221'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008501'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.15'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("D20").Value = '26.911.85'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.010'
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("D22").Value = '2.118.66'
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.34'
$ws.Range("E23").Value = '  -1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.372'
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.39'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.774'
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.087'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.38'
$ws.Range("E29").Value = '  +0.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.690'
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.725'
$ws.Range("E31").Value = '  +1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09129'
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05051'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7471'
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.981'
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.152'
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.227'
$ws.Range("E37").Value = '  +4.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.529'
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5594'
$ws.Range("E39").Value = '  +5.48%  '
$ws.Range("E40").Value = '  -1.94%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.619'
$ws.Range("E42").Value = '  +2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.69'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.554'
$ws.Range("E44").Value = '  +3.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1477'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4770'
$ws.Range("E46").Value = '  +3.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.05'
$ws.Range("E48").Value = '  +1.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.557'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.95'
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.00'
$ws.Range("E51").Value = '  -0.86%  '
